# Actualización automática 2025-08-29 10:50:09
$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": new sale of 367.8 in "PIEDRA SINTERIZADA" (col L) for row 18,
# and the running count in L23 moves from "2 de 21" to "3 de 21".
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L18").Value = 367.8
$wsGrupo.Range("L23").Value = "3 de 21"

# Sheet "VENTA MENSUAL": the same sale recorded in "agosto" (col F) for row 18,
# and the monthly total in F23 increases accordingly.
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F18").Value = 367.8
$wsMensual.Range("F23").Value = 1423.94

# Sheet "CUMPLIMIENTO MENSUAL": totals for group "OTROS" (row 2) and grand "TOTAL" (row 4)
# reflect the additional 367.8 in sales (VENTA / POR CUMPLIR / CUMPLIMIENTO).
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D2").Value = 74364.57000000001
$wsCumpl.Range("E2").Value = -74364.57000000001
$wsCumpl.Range("D4").Value = 77725.81000000001
$wsCumpl.Range("E4").Value = -62254.2507
$wsCumpl.Range("F4").Value = 5.023786451828421
